$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case Spanish connector words (de, del, la, las, el, los, y) in state/municipality names
$changes = @(
    @(8, "B", "Pabellón De Arteaga"),
    @(9, "B", "Rincón De Romos"),
    @(10, "B", "San Francisco De Los Romo"),
    @(11, "B", "San José De Gracia"),
    @(15, "B", "Playas De Rosarito"),
    @(32, "B", "Amatenango De La Frontera"),
    @(36, "B", "Bejucal De Ocampo"),
    @(38, "B", "Benemérito De Las Américas"),
    @(46, "B", "Chiapa De Corzo"),
    @(52, "B", "Comitán De Domínguez"),
    @(79, "B", "Marqués De Comillas"),
    @(80, "B", "Mazapa De Madero"),
    @(87, "B", "Ocozocoautla De Espinosa"),
    @(99, "B", "Salto De Agua"),
    @(101, "B", "San Cristóbal De Las Casas"),
    @(139, "B", "Guadalupe Y Calvo"),
    @(142, "B", "Hidalgo Del Parral"),
    @(154, "B", "San Francisco Del Oro"),
    @(156, "B", "Valle De Zaragoza"),
    @(158, "A", "Ciudad De México"),
    @(161, "B", "Cuajimalpa De Morelos"),
    @(176, "A", "Coahuila De Zaragoza"),
    @(190, "B", "San Juan De Sabinas"),
    @(203, "B", "Coneto De Comonfort"),
    @(215, "B", "Nombre De Dios"),
    @(219, "B", "Pánuco De Coronado"),
    @(227, "A", "Estado De México"),
    @(227, "B", "Acambay De Ruíz Castañeda"),
    @(230, "B", "Almoloya De Alquisiras"),
    @(231, "B", "Almoloya De Juárez"),
    @(236, "B", "Atizapán De Zaragoza"),
    @(241, "B", "Chapa De Mota"),
    @(244, "B", "Coacalco De Berriozábal"),
    @(250, "B", "Ecatepec De Morelos"),
    @(257, "B", "Ixtapan De La Sal"),
    @(258, "B", "Ixtapan Del Oro"),
    @(269, "B", "Naucalpan De Juárez"),
    @(274, "B", "San Felipe Del Progreso"),
    @(275, "B", "San José Del Rincón"),
    @(277, "B", "San Simón De Guerrero"),
    @(278, "B", "Soyaniquilpan De Juárez"),
    @(286, "B", "Tenango Del Valle"),
    @(294, "B", "Tlalnepantla De Baz"),
    @(299, "B", "Valle De Bravo"),
    @(300, "B", "Valle De Chalco Solidaridad"),
    @(303, "B", "Villa De Allende"),
    @(304, "B", "Villa Del Carbón"),
    @(315, "B", "Apaseo El Alto"),
    @(316, "B", "Apaseo El Grande"),
    @(324, "B", "Dolores Hidalgo Cuna De La Independencia Nacional"),
    @(328, "B", "Jaral Del Progreso"),
    @(335, "B", "Purísima Del Rincón"),
    @(340, "B", "San Diego De La Unión"),
    @(342, "B", "San Francisco Del Rincón"),
    @(344, "B", "San Luis De La Paz"),
    @(345, "B", "San Miguel De Allende"),
    @(346, "B", "Santa Cruz De Juventino Rosas"),
    @(347, "B", "Silao De La Victoria"),
    @(352, "B", "Valle De Santiago"),
    @(358, "B", "Acapulco De Juárez"),
    @(361, "B", "Ajuchitlán Del Progreso"),
    @(362, "B", "Alcozauca De Guerrero"),
    @(365, "B", "Atenango Del Río"),
    @(366, "B", "Atlamajalcingo Del Monte"),
    @(368, "B", "Atoyac De Álvarez"),
    @(369, "B", "Ayutla De Los Libres"),
    @(371, "B", "Buenavista De Cuéllar"),
    @(372, "B", "Chilapa De Álvarez"),
    @(373, "B", "Chilpancingo De Los Bravo"),
    @(374, "B", "Cochoapa El Grande"),
    @(377, "B", "Coyuca De Benítez"),
    @(378, "B", "Coyuca De Catalán"),
    @(381, "B", "Cuetzala Del Progreso"),
    @(382, "B", "Cutzamala De Pinzón"),
    @(389, "B", "Huitzuco De Los Figueroa"),
    @(390, "B", "Iguala De La Independencia"),
    @(392, "B", "Ixcateopan De Cuauhtémoc"),
    @(395, "B", "La Unión De Isidoro Montes De Oca"),
    @(400, "B", "Mártir De Cuilapan"),
    @(410, "B", "Taxco De Alarcón"),
    @(413, "B", "Tepecoacuilco De Trujano"),
    @(415, "B", "Tixtla De Guerrero"),
    @(419, "B", "Tlapa De Comonfort"),
    @(421, "B", "Técpan De Galeana"),
    @(426, "B", "Zihuatanejo De Azueta"),
    @(433, "B", "Agua Blanca De Iturbide"),
    @(437, "B", "Atotonilco El Grande"),
    @(443, "B", "Cuautepec De Hinojosa"),
    @(446, "B", "Huasca De Ocampo"),
    @(450, "B", "Huejutla De Reyes"),
    @(453, "B", "Jacala De Ledezma"),
    @(459, "B", "Mineral De La Reforma"),
    @(460, "B", "Mineral Del Chico"),
    @(461, "B", "Mineral Del Monte"),
    @(462, "B", "Mixquiahuala De Juárez"),
    @(463, "B", "Molango De Escamilla"),
    @(465, "B", "Nopala De Villagrán"),
    @(466, "B", "Omitlán De Juárez"),
    @(467, "B", "Pachuca De Soto"),
    @(469, "B", "Progreso De Obregón"),
    @(475, "B", "Santiago Tulantepec De Lugo Guerrero"),
    @(476, "B", "Santiago De Anaya"),
    @(480, "B", "Tenango De Doria"),
    @(482, "B", "Tepehuacán De Guerrero"),
    @(483, "B", "Tepeji Del Río De Ocampo"),
    @(484, "B", "Tezontepec De Aldama"),
    @(491, "B", "Tula De Allende"),
    @(492, "B", "Tulancingo De Bravo"),
    @(495, "B", "Zacualtipán De Ángeles"),
    @(500, "B", "Acatlán De Juárez"),
    @(501, "B", "Ahualulco De Mercado"),
    @(503, "B", "Atotonilco El Alto"),
    @(504, "B", "Autlán De Navarro"),
    @(508, "B", "Cañadas De Obregón"),
    @(515, "B", "Encarnación De Díaz"),
    @(517, "B", "Ixtlahuacán De Los Membrillos"),
    @(520, "B", "Jilotlán De Los Dolores"),
    @(523, "B", "Lagos De Moreno"),
    @(526, "B", "Ojuelos De Jalisco"),
    @(530, "B", "San Miguel El Alto"),
    @(534, "B", "Tamazula De Gordiano"),
    @(537, "B", "Tepatitlán De Morelos"),
    @(540, "B", "Tizapán El Alto"),
    @(541, "B", "Tlajomulco De Zúñiga"),
    @(547, "B", "Unión De San Antonio"),
    @(548, "B", "Unión De Tula"),
    @(551, "B", "Yahualica De González Gallo"),
    @(552, "B", "Zacoalco De Torres"),
    @(556, "B", "Zapotlán El Grande"),
    @(558, "A", "Michoacán De Ocampo"),
    @(573, "B", "Coalcomán De Vázquez Pallares"),
    @(629, "B", "Tiquicheo De Nicolás Romero"),
    @(651, "B", "Coatlán Del Río"),
    @(661, "B", "Puente De Ixtla"),
    @(666, "B", "Tetela Del Volcán"),
    @(667, "B", "Tlaltizapán De Zapata"),
    @(684, "B", "Santa María Del Oro"),
    @(692, "B", "Chiquihuitlán De Benito Juárez"),
    @(700, "B", "Mier Y Noriega"),
    @(705, "B", "San Nicolás De Los Garza"),
    @(710, "B", "Acatlán De Pérez Figueroa"),
    @(714, "B", "Chalcatongo De Hidalgo"),
    @(715, "B", "Chiquihuitlán De Benito Juárez"),
    @(718, "B", "Coicoyán De Las Flores"),
    @(719, "B", "Constancia Del Rosario"),
    @(722, "B", "Cuilápam De Guerrero"),
    @(724, "B", "El Barrio De La Soledad"),
    @(726, "B", "Guadalupe De Ramírez"),
    @(727, "B", "Guevea De Humboldt"),
    @(728, "B", "Heroica Ciudad De Ejutla De Crespo"),
    @(729, "B", "Heroica Ciudad De Huajuapan De León"),
    @(730, "B", "Heroica Ciudad De Juchitán De Zaragoza"),
    @(731, "B", "Heroica Ciudad De Tlaxiaco"),
    @(732, "B", "Huajuapan De León"),
    @(733, "B", "Huautla De Jiménez"),
    @(735, "B", "Ixtlán De Juárez"),
    @(739, "B", "Mazatlán Villa De Flores"),
    @(740, "B", "Miahuatlán De Porfirio Díaz"),
    @(742, "B", "Mártires De Tacubaya"),
    @(744, "B", "Nejapa De Madero"),
    @(745, "B", "Oaxaca De Juárez"),
    @(746, "B", "Ocotlán De Morelos"),
    @(747, "B", "Pinotepa De Don Luis"),
    @(749, "B", "Putla Villa De Guerrero"),
    @(750, "B", "Reforma De Pineda"),
    @(757, "B", "San Antonino El Alto"),
    @(761, "B", "San Baltazar Yatzachi El Bajo"),
    @(766, "B", "San Felipe Jalapa De Díaz"),
    @(770, "B", "San Francisco Del Mar"),
    @(778, "B", "San José Del Progreso"),
    @(780, "B", "San Juan Bautista Lo De Soto"),
    @(798, "B", "San Juan De Los Cués"),
    @(820, "B", "San Miguel Del Puerto"),
    @(821, "B", "San Miguel El Grande"),
    @(833, "B", "San Pedro El Alto"),
    @(863, "B", "Santa María Jalapa Del Marqués"),
    @(874, "B", "Santa María Del Tule"),
    @(904, "B", "Santo Domingo De Morelos"),
    @(908, "B", "Tataltepec De Valdés"),
    @(909, "B", "Teococuilco De Marcos Pérez"),
    @(910, "B", "Teotitlán De Flores Magón"),
    @(911, "B", "Tlacolula De Matamoros"),
    @(913, "B", "Villa Sola De Vega"),
    @(914, "B", "Villa De Etla"),
    @(915, "B", "Villa De Tututepec"),
    @(916, "B", "Villa De Tututepec De Melchor Ocampo"),
    @(917, "B", "Villa De Zaachila"),
    @(919, "B", "Zimatlán De Álvarez"),
    @(933, "B", "Chalchicomula De Sesma"),
    @(942, "B", "Cuetzalan Del Progreso"),
    @(954, "B", "Ixcamilpa De Guerrero"),
    @(955, "B", "Izúcar De Matamoros"),
    @(959, "B", "Los Reyes De Juárez"),
    @(960, "B", "Mazapiltepec De Juárez"),
    @(965, "B", "Palmar De Bravo"),
    @(975, "B", "San Nicolás De Los Ranchos"),
    @(977, "B", "San Salvador El Seco"),
    @(978, "B", "San Salvador El Verde"),
    @(983, "B", "Tepanco De López"),
    @(984, "B", "Tepango De Rodríguez"),
    @(985, "B", "Tepatlaxco De Hidalgo"),
    @(988, "B", "Tepexi De Rodríguez"),
    @(989, "B", "Tetela De Ocampo"),
    @(990, "B", "Teteles De Avila Castillo"),
    @(993, "B", "Tlacotepec De Benito Juárez"),
    @(1015, "B", "Amealco De Bonfil"),
    @(1017, "B", "Cadereyta De Montes"),
    @(1024, "B", "Jalpan De Serra"),
    @(1025, "B", "Landa De Matamoros"),
    @(1026, "B", "Pinal De Amoles"),
    @(1029, "B", "San Juan Del Río"),
    @(1041, "B", "Armadillo De Los Infante"),
    @(1042, "B", "Axtla De Terrazas"),
    @(1049, "B", "Ciudad Del Maíz"),
    @(1058, "B", "Mexquitic De Carmona"),
    @(1063, "B", "San Ciro De Acosta"),
    @(1069, "B", "Santa María Del Río"),
    @(1071, "B", "Soledad De Graciano Sánchez"),
    @(1078, "B", "Tanquián De Escobedo"),
    @(1083, "B", "Villa De Arista"),
    @(1084, "B", "Villa De Arriaga"),
    @(1085, "B", "Villa De Guadalupe"),
    @(1086, "B", "Villa De Ramos"),
    @(1087, "B", "Villa De Reyes"),
    @(1112, "B", "Nacozari De García"),
    @(1154, "B", "Soto La Marina"),
    @(1161, "B", "Apetatitlán De Antonio Carvajal"),
    @(1167, "B", "San Pablo Del Monte"),
    @(1174, "A", "Veracruz De Ignacio De La Llave"),
    @(1182, "B", "Amatlán De Los Reyes"),
    @(1191, "B", "Boca Del Río"),
    @(1193, "B", "Camarón De Tejeda"),
    @(1197, "B", "Cazones De Herrera"),
    @(1211, "B", "Cosamaloapan De Carpio"),
    @(1230, "B", "Hueyapan De Ocampo"),
    @(1231, "B", "Huiloapan De Cuauhtémoc"),
    @(1232, "B", "Ignacio De La Llave"),
    @(1235, "B", "Ixhuacán De Los Reyes"),
    @(1236, "B", "Ixhuatlán De Madero"),
    @(1237, "B", "Ixhuatlán Del Café"),
    @(1238, "B", "Ixhuatlán Del Sureste"),
    @(1252, "B", "Las Vigas De Ramírez"),
    @(1253, "B", "Lerdo De Tejada"),
    @(1258, "B", "Martínez De La Torre"),
    @(1264, "B", "Mixtla De Altamirano"),
    @(1266, "B", "Nanchital De Lázaro Cárdenas Del Río"),
    @(1276, "B", "Ozuluama De Mascareñas"),
    @(1279, "B", "Paso De Ovejas"),
    @(1280, "B", "Paso Del Macho"),
    @(1283, "B", "Poza Rica De Hidalgo"),
    @(1290, "B", "Sayula De Alemán"),
    @(1292, "B", "Soledad De Doblado"),
    @(1296, "B", "Tatahuicapan De Juárez"),
    @(1314, "B", "Tlacotepec De Mejía"),
    @(1324, "B", "Vega De Alatorre"),
    @(1332, "B", "Zontecomatlán De López Y Fuentes"),
    @(1345, "B", "Cañitas De Felipe Pescador"),
    @(1347, "B", "Concepción Del Oro"),
    @(1348, "B", "El Plateado De Joaquín Amaro"),
    @(1362, "B", "Moyahua De Estrada"),
    @(1363, "B", "Nochistlán De Mejía"),
    @(1364, "B", "Noria De Ángeles"),
    @(1372, "B", "Teúl De González Ortega"),
    @(1373, "B", "Tlaltenango De Sánchez Román"),
    @(1378, "B", "Villa De Cos")
)

foreach ($chg in $changes) {
    $r = $chg[0]
    $c = $chg[1]
    $v = $chg[2]
    $ws.Range("$c$r").Value = $v
}

# 3. Fix the TOTAL row label
$ws.Range("A1382").Value = "Total"

# 4. Remove trailing metadata/footer rows (1384-1388), which also shrinks the used range
$ws.Range("A1384:A1388").EntireRow.Delete()
